# Update example to new accnr format in excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A2 from the numeric accnr 202201037 to the new text-based
# accnr format "A2022/12345"
$ws.Range("A2").Value = "A2022/12345"

# Add a new row 3 that mirrors row 2 but uses the new accnr format
# without the slash: "A202212345"
$ws.Range("A3").Value = "A202212345"
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

$ws.Range("C13").Select()
